$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 10 ("Objetivos:"): value changes from the old "Aulas expositivas..." text
# to the first responsible-teacher identifier string.
$ws.Range("B10").Value = "7290967 - Emerson Gonçalves de Melo"
$ws.Range("C10").Value = "7290967 - Emerson Gonçalves de Melo"

# Row 13 ("Programa resumido:"): value changes from "90 h" to the activation
# date string "01/01/2023". Because that text looks like a date, typing it
# directly would make Excel auto-convert the cell to a date value/format.
# Cell B8 already holds this exact string as plain text, so we copy it from
# there and paste only the value, which keeps it as text without touching
# the cell's number format/style (exactly like genuine Excel behavior when
# duplicating an existing text cell).
$ws.Range("B8").Copy()
$ws.Range("B13").PasteSpecial(-4163)
$ws.Range("B8").Copy()
$ws.Range("C13").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Row 15 ("Programa:"): value changes from "Semestral" to the first teacher
# identifier string (same one used in row 10).
$ws.Range("B15").Value = "7290967 - Emerson Gonçalves de Melo"
$ws.Range("C15").Value = "7290967 - Emerson Gonçalves de Melo"

# Row 18 ("Método:"): value changes from "01/01/2023" to the second teacher
# identifier string.
$ws.Range("B18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
